$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 999.6667
$ws.Range("I19").Value = 998.5
$ws.Range("J19").Value = 1002
$ws.Range("K19").Value = 998.5
$ws.Range("L19").Value = 1002
$ws.Range("M19").Value = -823.5
$ws.Range("N19").Value = -1352
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 487.75
$ws.Range("I33").Value = 577
$ws.Range("K33").Value = 577
$ws.Range("M33").Value = -348
# Row 34 (Leve Item ID 2160)
$ws.Range("H34").Value = 6035.625
$ws.Range("J34").Value = 10998.75
$ws.Range("L34").Value = 10998.75
$ws.Range("N34").Value = -11404.75
# Row 36 (Leve Item ID 2160)
$ws.Range("H36").Value = 6035.625
$ws.Range("J36").Value = 10998.75
$ws.Range("L36").Value = 10998.75
$ws.Range("N36").Value = -12428.75
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 5005.857
$ws.Range("I40").Value = 3057.6667
$ws.Range("J40").Value = 6467
$ws.Range("K40").Value = 3057.6667
$ws.Range("L40").Value = 6467
$ws.Range("M40").Value = -2882.6667
$ws.Range("N40").Value = -6817
# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 5271.6665
$ws.Range("I86").Value = 4408.7144
$ws.Range("J86").Value = 6479.8
$ws.Range("K86").Value = 4408.7144
$ws.Range("L86").Value = 6479.8
$ws.Range("M86").Value = -3285.7144
$ws.Range("N86").Value = -8725.799999999999
# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 5271.6665
$ws.Range("I89").Value = 4408.7144
$ws.Range("J89").Value = 6479.8
$ws.Range("K89").Value = 22043.572
$ws.Range("L89").Value = 32399
$ws.Range("M89").Value = -16427.572
$ws.Range("N89").Value = -43631
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 13910.546
$ws.Range("J116").Value = 15002.4
$ws.Range("L116").Value = 15002.4
$ws.Range("N116").Value = -21886.4
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4042.6956
$ws.Range("I138").Value = 5839.2
$ws.Range("K138").Value = 17517.6
$ws.Range("M138").Value = -12377.6

$ws = $wb.Worksheets.Item("ARM")
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 4487
$ws.Range("J110").Value = 6253
$ws.Range("L110").Value = 6253
$ws.Range("N110").Value = -10343
# Row 124 (Leve Item ID 34252)
$ws.Range("H124").Value = 46187.5
$ws.Range("J124").Value = 46187.5
$ws.Range("L124").Value = 46187.5
$ws.Range("N124").Value = -56007.5
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4351.8
$ws.Range("I132").Value = 2711.7273
$ws.Range("K132").Value = 8135.1819
$ws.Range("M132").Value = -5605.1819

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3859
$ws.Range("I134").Value = 1999.8334
$ws.Range("K134").Value = 5999.5002
$ws.Range("M134").Value = -3464.5002

$ws = $wb.Worksheets.Item("CRP")
# Row 10 (Leve Item ID 1997)
$ws.Range("H10").Value = 1353.375
$ws.Range("J10").Value = 1635.3334
$ws.Range("L10").Value = 1635.3334
$ws.Range("N10").Value = -1913.3334
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 4021.4443
$ws.Range("J16").Value = 5333.3335
$ws.Range("L16").Value = 5333.3335
$ws.Range("N16").Value = -5907.3335
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 24066.451
$ws.Range("I31").Value = 3261.805
$ws.Range("J31").Value = 109365.5
$ws.Range("K31").Value = 3261.805
$ws.Range("L31").Value = 109365.5
$ws.Range("M31").Value = -2966.805
$ws.Range("N31").Value = -109955.5
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 24066.451
$ws.Range("I34").Value = 3261.805
$ws.Range("J34").Value = 109365.5
$ws.Range("K34").Value = 3261.805
$ws.Range("L34").Value = 109365.5
$ws.Range("M34").Value = -3059.805
$ws.Range("N34").Value = -109769.5
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1581.3572
$ws.Range("I107").Value = 2149.1667
$ws.Range("J107").Value = 1155.5
$ws.Range("K107").Value = 2149.1667
$ws.Range("L107").Value = 1155.5
$ws.Range("M107").Value = -229.1667000000002
$ws.Range("N107").Value = -4995.5
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 4021.4443
$ws.Range("J113").Value = 5333.3335
$ws.Range("L113").Value = 5333.3335
$ws.Range("N113").Value = -9673.333500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 64 (Leve Item ID 12861)
$ws.Range("H64").Value = 166672980
$ws.Range("I64").Value = 333336100
$ws.Range("J64").Value = 9854.666999999999
$ws.Range("K64").Value = 1000008300
$ws.Range("L64").Value = 29564.001
$ws.Range("M64").Value = -1000008030
$ws.Range("N64").Value = -30104.001
# Row 67 (Leve Item ID 12861)
$ws.Range("H67").Value = 166672980
$ws.Range("I67").Value = 333336100
$ws.Range("J67").Value = 9854.666999999999
$ws.Range("K67").Value = 1000008300
$ws.Range("L67").Value = 29564.001
$ws.Range("M67").Value = -1000007364
$ws.Range("N67").Value = -31436.001
# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 6679
$ws.Range("J69").Value = 10014
$ws.Range("L69").Value = 30042
$ws.Range("N69").Value = -31664
# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 6679
$ws.Range("J72").Value = 10014
$ws.Range("L72").Value = 90126
$ws.Range("N72").Value = -98238
# Row 81 (Leve Item ID 12843)
$ws.Range("H81").Value = 3058.3333
$ws.Range("I81").Value = 1700
$ws.Range("K81").Value = 5100
$ws.Range("M81").Value = -3977
# Row 84 (Leve Item ID 12843)
$ws.Range("H84").Value = 3058.3333
$ws.Range("I84").Value = 1700
$ws.Range("K84").Value = 15300
$ws.Range("M84").Value = -9684
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 1062.1666
$ws.Range("J113").Value = 975.7778
$ws.Range("L113").Value = 2927.3334
$ws.Range("N113").Value = -7267.3334
# Row 136 (Leve Item ID 44093)
$ws.Range("H136").Value = 1957.25
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 7022.857
$ws.Range("I137").Value = 2014.5
$ws.Range("J137").Value = 9026.200000000001
$ws.Range("K137").Value = 6043.5
$ws.Range("L137").Value = 27078.6
$ws.Range("M137").Value = -943.5
$ws.Range("N137").Value = -37278.60000000001
# Row 139 (Leve Item ID 44102)
$ws.Range("H139").Value = 4199.4707
$ws.Range("I139").Value = 1646.5834
$ws.Range("J139").Value = 10326.4
$ws.Range("K139").Value = 4939.7502
$ws.Range("L139").Value = 30979.2
$ws.Range("M139").Value = 200.2497999999996
$ws.Range("N139").Value = -41259.2

$ws = $wb.Worksheets.Item("GSM")
# Row 98 (Leve Item ID 18359)
$ws.Range("H98").Value = 19333.334
$ws.Range("J98").Value = 19333.334
$ws.Range("L98").Value = 19333.334
$ws.Range("N98").Value = -25323.334
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 4990.278
$ws.Range("J113").Value = 5917.6
$ws.Range("L113").Value = 5917.6
$ws.Range("N113").Value = -10257.6

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 7051.5415
$ws.Range("I40").Value = 6011.4
$ws.Range("J40").Value = 12252.25
$ws.Range("K40").Value = 6011.4
$ws.Range("L40").Value = 12252.25
$ws.Range("M40").Value = -5875.4
$ws.Range("N40").Value = -12524.25
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 4830.1
$ws.Range("I46").Value = 1900
$ws.Range("K46").Value = 1900
$ws.Range("M46").Value = -1712
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 3377.1333
$ws.Range("J61").Value = 13164.667
$ws.Range("L61").Value = 13164.667
$ws.Range("N61").Value = -13568.667
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3553.8572
$ws.Range("J68").Value = 3930.6667
$ws.Range("L68").Value = 3930.6667
$ws.Range("N68").Value = -5428.6667
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3553.8572
$ws.Range("J71").Value = 3930.6667
$ws.Range("L71").Value = 19653.3335
$ws.Range("N71").Value = -27141.3335
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 3377.1333
$ws.Range("J113").Value = 13164.667
$ws.Range("L113").Value = 13164.667
$ws.Range("N113").Value = -17504.667

$ws = $wb.Worksheets.Item("WVR")
# Row 70 (Leve Item ID 11979)
$ws.Range("H70").Value = 21582.834
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73 (Leve Item ID 11979)
$ws.Range("H73").Value = 21582.834
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

